$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2749.5
$ws.Range("I40").Value = 499
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 499
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -324
$ws.Range("N40").Value = -5350
$ws.Range("H123").Value = 180000
$ws.Range("J123").Value = 180000
$ws.Range("L123").Value = 180000
$ws.Range("N123").Value = -189800
$ws.Range("H129").Value = 2863.3
$ws.Range("J129").Value = 2991.7144
$ws.Range("L129").Value = 8975.143199999999
$ws.Range("N129").Value = -18975.1432
$ws.Range("H132").Value = 3198.5
$ws.Range("I132").Value = 2442.8333
$ws.Range("K132").Value = 7328.499899999999
$ws.Range("M132").Value = -4798.499899999999
$ws.Range("H135").Value = 2110.6924
$ws.Range("I135").Value = 2341
$ws.Range("J135").Value = 1343
$ws.Range("K135").Value = 21069
$ws.Range("L135").Value = 12087
$ws.Range("M135").Value = -18534
$ws.Range("N135").Value = -17157
$ws.Range("H137").Value = 2427.9756
$ws.Range("I137").Value = 1615.1111
$ws.Range("J137").Value = 3995.6428
$ws.Range("K137").Value = 4845.3333
$ws.Range("L137").Value = 11986.9284
$ws.Range("M137").Value = -2295.3333
$ws.Range("N137").Value = -17086.9284
$ws.Range("H138").Value = 4794.737
$ws.Range("J138").Value = 5274.7334
$ws.Range("L138").Value = 15824.2002
$ws.Range("N138").Value = -26104.2002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 1999.5
$ws.Range("I11").Value = 999
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 999
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -855
$ws.Range("N11").Value = -3288
$ws.Range("H61").Value = 2199.3333
$ws.Range("I61").Value = 1962.25
$ws.Range("K61").Value = 1962.25
$ws.Range("M61").Value = -1750.25
$ws.Range("H102").Value = 1499
$ws.Range("I102").Value = 1499
$ws.Range("K102").Value = 1499
$ws.Range("M102").Value = 123
$ws.Range("H132").Value = 2109.913
$ws.Range("I132").Value = 1315.4
$ws.Range("K132").Value = 3946.2
$ws.Range("M132").Value = -1416.2
$ws.Range("H136").Value = 2199.3333
$ws.Range("I136").Value = 1962.25
$ws.Range("K136").Value = 5886.75
$ws.Range("M136").Value = -3336.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 582.5
$ws.Range("I80").Value = 1214.5
$ws.Range("J80").Value = 329.7
$ws.Range("K80").Value = 1214.5
$ws.Range("L80").Value = 329.7
$ws.Range("M80").Value = -216.5
$ws.Range("N80").Value = -2325.7
$ws.Range("H83").Value = 582.5
$ws.Range("I83").Value = 1214.5
$ws.Range("J83").Value = 329.7
$ws.Range("K83").Value = 6072.5
$ws.Range("L83").Value = 1648.5
$ws.Range("M83").Value = -1080.5
$ws.Range("N83").Value = -11632.5
$ws.Range("H134").Value = 4562.476
$ws.Range("I134").Value = 4411.2104
$ws.Range("K134").Value = 13233.6312
$ws.Range("M134").Value = -10698.6312
$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200
$ws.Range("H138").Value = 106292
$ws.Range("J138").Value = 106292
$ws.Range("L138").Value = 106292
$ws.Range("N138").Value = -116572

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2580.9546
$ws.Range("I31").Value = 2613.4285
$ws.Range("J31").Value = 1899
$ws.Range("K31").Value = 2613.4285
$ws.Range("L31").Value = 1899
$ws.Range("M31").Value = -2318.4285
$ws.Range("N31").Value = -2489
$ws.Range("H34").Value = 2580.9546
$ws.Range("I34").Value = 2613.4285
$ws.Range("J34").Value = 1899
$ws.Range("K34").Value = 2613.4285
$ws.Range("L34").Value = 1899
$ws.Range("M34").Value = -2411.4285
$ws.Range("N34").Value = -2303
$ws.Range("H52").Value = 133300
$ws.Range("J52").Value = 139950
$ws.Range("L52").Value = 139950
$ws.Range("N52").Value = -140538
$ws.Range("H58").Value = 2748.75
$ws.Range("I58").Value = 2748.75
$ws.Range("K58").Value = 2748.75
$ws.Range("M58").Value = -2545.75
$ws.Range("H94").Value = 1261.5
$ws.Range("J94").Value = 1111
$ws.Range("L94").Value = 1111
$ws.Range("N94").Value = -2013
$ws.Range("H105").Value = 2636
$ws.Range("I105").Value = 2345.818
$ws.Range("K105").Value = 2345.818
$ws.Range("M105").Value = -598.8180000000002
$ws.Range("H107").Value = 1680.3889
$ws.Range("I107").Value = 1018.2727
$ws.Range("K107").Value = 1018.2727
$ws.Range("M107").Value = 901.7273
$ws.Range("H136").Value = 2748.75
$ws.Range("I136").Value = 2748.75
$ws.Range("K136").Value = 8246.25
$ws.Range("M136").Value = -5696.25
$ws.Range("H141").Value = 110400.664
$ws.Range("J141").Value = 112480.8
$ws.Range("L141").Value = 112480.8
$ws.Range("N141").Value = -122840.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 900
$ws.Range("I13").Value = 900
$ws.Range("K13").Value = 2700
$ws.Range("M13").Value = -2532
$ws.Range("H69").Value = 1800
$ws.Range("I69").Value = 1600
$ws.Range("K69").Value = 4800
$ws.Range("M69").Value = -3989
$ws.Range("H72").Value = 1800
$ws.Range("I72").Value = 1600
$ws.Range("K72").Value = 14400
$ws.Range("M72").Value = -10344
$ws.Range("H80").Value = 999999
$ws.Range("J80").Value = 999999
$ws.Range("L80").Value = 2999997
$ws.Range("N80").Value = -3001869
$ws.Range("H83").Value = 999999
$ws.Range("J83").Value = 999999
$ws.Range("L83").Value = 8999991
$ws.Range("N83").Value = -9009351
$ws.Range("H131").Value = 2087.0557
$ws.Range("I131").Value = 1599
$ws.Range("J131").Value = 2115.7646
$ws.Range("K131").Value = 4797
$ws.Range("L131").Value = 6347.293799999999
$ws.Range("M131").Value = 243
$ws.Range("N131").Value = -16427.2938

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2670.3333
$ws.Range("I132").Value = 2074.5
$ws.Range("K132").Value = 6223.5
$ws.Range("M132").Value = -3693.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840
$ws.Range("H132").Value = 4157.6665
$ws.Range("I132").Value = 3570.1667
$ws.Range("K132").Value = 10710.5001
$ws.Range("M132").Value = -8180.500100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("H132").Value = 2186.842
$ws.Range("I132").Value = 1659.625
$ws.Range("J132").Value = 4998.6665
$ws.Range("K132").Value = 4978.875
$ws.Range("L132").Value = 14995.9995
$ws.Range("M132").Value = -2448.875
$ws.Range("N132").Value = -20055.9995
$ws.Range("H137").Value = 87599.75
$ws.Range("J137").Value = 87599.75
$ws.Range("L137").Value = 87599.75
$ws.Range("N137").Value = -97799.75
$ws.Range("M11").ClearContents()
